$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.213.79"
$ws.Range("E2").Value = "  -0.28%  "

# Row 3
$ws.Range("D3").Value = "2.642.36"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.08%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.92%  "

# Row 9
$ws.Range("E9").Value = "  -1.82%  "

# Row 10
$ws.Range("E10").Value = "  -1.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "

# Row 12
$ws.Range("E12").Value = "  -1.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.31%  "

# Row 14
$ws.Range("D14").Value = "3.124.94"
$ws.Range("E14").Value = "  -0.02%  "

# Row 15
$ws.Range("E15").Value = "  -2.79%  "

# Row 16
$ws.Range("D16").Value = "68.162.98"
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("D17").Value = "2.650.08"
$ws.Range("E17").Value = "  -0.48%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "358.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.05%  "

# Row 20
$ws.Range("E20").Value = "  -2.61%  "

# Row 21
$ws.Range("E21").Value = "  -0.22%  "

# Row 22
$ws.Range("E22").Value = "  -3.77%  "

# Row 23
$ws.Range("E23").Value = "  -0.81%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "

# Row 25
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("E26").Value = "  -1.45%  "

# Row 27
$ws.Range("D27").Value = "2.799.07"
$ws.Range("E27").Value = "  +0.88%  "

# Row 28
$ws.Range("E28").Value = "  -3.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "556.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.66%  "

# Row 32
$ws.Range("E32").Value = "  -3.68%  "

# Row 33
$ws.Range("E33").Value = "  +0.02%  "

# Row 34
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.73%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.127"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.68%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.29%  "

# Row 38
$ws.Range("E38").Value = "  +1.18%  "

# Row 39
$ws.Range("E39").Value = "  -1.68%  "

# Row 40
$ws.Range("E40").Value = "  -3.03%  "

# Row 41
$ws.Range("E41").Value = "  -2.19%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.03%  "

# Row 43
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0317"
$ws.Range("E44").Value = "  -6.99%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "156.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "

# Row 46
$ws.Range("E46").Value = "  -0.22%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "

# Row 48
$ws.Range("E48").Value = "  -3.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0772"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.17%  "

# Row 50
$ws.Range("E50").Value = "  -1.03%  "

# Row 51
$ws.Range("E51").Value = "  -0.48%  "
